$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 937 (mirrors header pattern repeated at rows 1,37,73,...,901)
$headers = @("States/UT", "Active Cases", "Active Cases Since Yesterday", "Recovered Cases", "Recovered Cases Since Yesterday", "Deceased Cases", "Deceased Cases Since Yesterday", "Date")
for ($c = 1; $c -le 8; $c++) {
    $ws.Cells.Item(937, $c).Value = $headers[$c - 1]
}

# Copy the header formatting (bold, thin border, center/top alignment) from the
# previous header block (row 901) so the new header row 937 matches the
# established repeating style used throughout the sheet.
$ws.Range("A901:H901").Copy() | Out-Null
$ws.Range("A937:H937").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data rows 938-972: one row per state/UT for 13-10-2020
$data = @(
    @('Andaman and Nicobar Islands', 198, 12, 3770, 6, 55, 0, '13-10-2020'),
    @('Andhra Pradesh', 43983, -2312, 708712, 5504, 6256, 32, '13-10-2020'),
    @('Arunachal Pradesh', 2940, 49, 9403, 171, 24, 0, '13-10-2020'),
    @('Assam', 28439, 54, 166039, 1457, 826, 10, '13-10-2020'),
    @('Bihar', 10669, -375, 185911, 1370, 955, 9, '13-10-2020'),
    @('Chandigarh', 1170, -14, 11898, 111, 192, 1, '13-10-2020'),
    @('Chhattisgarh', 27421, 73, 116540, 2769, 1286, 33, '13-10-2020'),
    @('Dadra and Nagar Haveli and Daman and Diu', 102, 0, 3064, 5, 2, 0, '13-10-2020'),
    @('Delhi', 20535, -1166, 284844, 2975, 5809, 40, '13-10-2020'),
    @('Goa', 4465, -191, 33698, 495, 511, 4, '13-10-2020'),
    @('Gujarat', 15414, -281, 133615, 1442, 3574, 8, '13-10-2020'),
    @('Haryana', 10401, -172, 131228, 1225, 1592, 13, '13-10-2020'),
    @('Himachal Pradesh', 2637, -50, 14690, 219, 251, 1, '13-10-2020'),
    @('Jammu and Kashmir', 9992, -474, 72706, 861, 1333, 11, '13-10-2020'),
    @('Jharkhand', 7776, -391, 84461, 890, 798, 11, '13-10-2020'),
    @('Karnataka', 115795, -4494, 592084, 12030, 10036, 70, '13-10-2020'),
    @('Kerala', 94473, -1928, 199634, 7836, 1025, 22, '13-10-2020'),
    @('Ladakh', 961, -19, 4126, 89, 64, 0, '13-10-2020'),
    @('Madhya Pradesh', 14932, -245, 130721, 1702, 2645, 21, '13-10-2020'),
    @('Maharashtra', 212905, -8732, 1281896, 15656, 40514, 165, '13-10-2020'),
    @('Manipur', 2756, 25, 10707, 203, 93, 2, '13-10-2020'),
    @('Meghalaya', 2434, -44, 5273, 131, 64, 1, '13-10-2020'),
    @('Mizoram', 156, -18, 2046, 36, 0, 0, '13-10-2020'),
    @('Nagaland', 1409, 150, 5813, 70, 18, 1, '13-10-2020'),
    @('Odisha', 23430, -172, 230192, 2577, 1040, 18, '13-10-2020'),
    @('Puducherry', 4617, -78, 26555, 264, 565, 2, '13-10-2020'),
    @('Punjab', 8576, -699, 112099, 1234, 3860, 27, '13-10-2020'),
    @('Rajasthan', 21671, 259, 137848, 1858, 1665, 15, '13-10-2020'),
    @('Sikkim', 384, 0, 2925, 5, 57, 2, '13-10-2020'),
    @('Tamil Nadu', 43747, -348, 607203, 5165, 10314, 62, '13-10-2020'),
    @('Telengana', 24208, -306, 189351, 2009, 1233, 5, '13-10-2020'),
    @('Tripura', 3738, -4, 24623, 220, 317, 1, '13-10-2020'),
    @('Uttarakhand', 6976, -397, 47609, 678, 762, 15, '13-10-2020'),
    @('Uttar Pradesh', 38815, -1204, 393908, 3342, 6438, 44, '13-10-2020'),
    @('West Bengal', 30604, 368, 262103, 3155, 5682, 60, '13-10-2020')
)

$r = 938
foreach ($row in $data) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

Write-Host "Added rows 937-972 (header + 35 states) for 13-10-2020"
